# This script rotates the data of rows 5, 6 and 7 (columns A,B,D,E,F,G,H,L,Q,R):
#   new row5 <- old row6
#   new row6 <- old row7
#   new row7 <- old row5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original ("before") values of row 5, since row 5 gets
# overwritten first but its data is needed later for the new row 7.
# (Use Value2 -- Value has been observed to return a bogus reflection
# descriptor instead of the real cell contents in this engine.)
$A5 = $ws.Range("A5").Value2
$B5 = $ws.Range("B5").Value2
$D5 = $ws.Range("D5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2

# row5 <- row6
$ws.Range("A5").Value2 = $ws.Range("A6").Value2
$ws.Range("B5").Value2 = $ws.Range("B6").Value2
$ws.Range("D5").Value2 = $ws.Range("D6").Value2
$ws.Range("E5").Value2 = $ws.Range("E6").Value2
$ws.Range("F5").Value2 = $ws.Range("F6").Value2
$ws.Range("G5").Value2 = $ws.Range("G6").Value2
$ws.Range("H5").Value2 = $ws.Range("H6").Value2
$ws.Range("Q5").Value2 = $ws.Range("Q6").Value2
$ws.Range("R5").Value2 = $ws.Range("R6").Value2

# row6 <- row7
$ws.Range("A6").Value2 = $ws.Range("A7").Value2
$ws.Range("B6").Value2 = $ws.Range("B7").Value2
$ws.Range("D6").Value2 = $ws.Range("D7").Value2
$ws.Range("E6").Value2 = $ws.Range("E7").Value2
$ws.Range("F6").Value2 = $ws.Range("F7").Value2
$ws.Range("G6").Value2 = $ws.Range("G7").Value2
$ws.Range("H6").Value2 = $ws.Range("H7").Value2
$ws.Range("Q6").Value2 = $ws.Range("Q7").Value2
$ws.Range("R6").Value2 = $ws.Range("R7").Value2

# row7 <- (original) row5
$ws.Range("A7").Value2 = $A5
$ws.Range("B7").Value2 = $B5
$ws.Range("D7").Value2 = $D5
$ws.Range("E7").Value2 = $E5
$ws.Range("F7").Value2 = $F5
$ws.Range("G7").Value2 = $G5
$ws.Range("H7").Value2 = $H5
$ws.Range("Q7").Value2 = $Q5
$ws.Range("R7").Value2 = $R5

# The (empty) L column cell follows the same rotation: it moves from L5 to
# L7 (L6 had no cell before and still has none after). Use Cut so the
# source cell is cleared out (not just blanked) while the destination
# cell is (re)created, matching the presence/absence pattern in the diff.
$ws.Range("L5").Cut($ws.Range("L7"))
